# Added the winch: a new "winchMotor" PWM output on the Outputs sheet and a
# new "winchSwitch" digital input on the Other Inputs sheet. Also renamed the
# three sheets to their descriptive names (Outputs / Joysticks / Other Inputs).

$wb = $excel.ActiveWorkbook

$wsOutputs = $wb.Worksheets.Item(1)
$wsJoysticks = $wb.Worksheets.Item(2)
$wsOther = $wb.Worksheets.Item(3)

# --- Rename the sheets -----------------------------------------------------
$wsOutputs.Name = "Outputs"
$wsJoysticks.Name = "Joysticks"
$wsOther.Name = "Other Inputs"

# --- Outputs sheet: insert the winchMotor row right above "solenoids[]" ----
# (row 7). Row 8 above (row 6, "ballRoller") doesn't carry the border-style
# that the solenoid header row does, so pull a known-good style (row 6,
# columns A-D are style "7") onto the new row first, then fill in values.
$wsOutputs.Rows.Item(7).Insert()

$wsOutputs.Range("A6:E6").Copy()
$wsOutputs.Range("A7").PasteSpecial(-4104)

$wsOutputs.Range("A7").Value = "winchMotor"
$wsOutputs.Range("B7").Value = "PWM"
$wsOutputs.Range("C7").Value = "???"
$wsOutputs.Range("D7").Value = "???"
$wsOutputs.Range("E7").Value = "Lift us up"

# --- Other Inputs sheet: insert the winchSwitch row above "stoppedModeS" ---
# (row 8).
$wsOther.Rows.Item(8).Insert()

$wsOther.Range("B7").Copy()
$wsOther.Range("E8").PasteSpecial(-4104)

$wsOther.Range("A8").Value = "winchSwitch"
$wsOther.Range("B8").Value = "DigitalInput"
$wsOther.Range("C8").Value = "???"
$wsOther.Range("D8").Value = "???"
$wsOther.Range("E8").Value = "Do we want to turn on the winch motor?"

# --- Selections (cosmetic, matches the saved workbook view state) ---------
$wsOutputs.Range("B20").Select()
$wsOther.Range("E8").Select()
